# ------------------------------------------------------------------
# "changed how I define z score in the network plot to make it more
#  reasonable" -- adds a 4th patient sheet (Pt4), renames the
#  "Stim seizure elecs" header to "Clinical Stim seizure elecs" on
#  every sheet, and fixes up a couple of font-style inconsistencies
#  on the existing sheets (F1/M1 header cells, a couple of data
#  cells on Pt2).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

$ws1 = $wb.Worksheets.Item("Pt1")
$ws2 = $wb.Worksheets.Item("Pt2")
$ws3 = $wb.Worksheets.Item("Pt3")

# --------------------------------------------------------------
# 1. Rename the shared header string everywhere it is used.
# --------------------------------------------------------------
$ws1.Range("F1").Value = "Clinical Stim seizure elecs"
$ws2.Range("F1").Value = "Clinical Stim seizure elecs"
$ws3.Range("F1").Value = "Clinical Stim seizure elecs"

# --------------------------------------------------------------
# 2. Fix up the font-style swaps on the existing sheets.
#
#    F1 goes from the "Arial/theme-1" look to the "default" look,
#    and M1 (plus a few Pt2 data cells) goes the other way. Grab a
#    style-2 donor cell (Pt1!M1, before we touch it) and a style-1
#    donor cell (Pt1!A1, never touched) and stamp the right format
#    onto every target with Copy + PasteSpecial(xlPasteFormats) so
#    the underlying font table stays deduplicated exactly like
#    Excel would produce.
# --------------------------------------------------------------
$style2Donor = $ws1.Range("M1")
$style1Donor = $ws1.Range("A1")

$style2Donor.Copy()
$ws1.Range("F1").PasteSpecial($xlPasteFormats)
$ws2.Range("F1").PasteSpecial($xlPasteFormats)
$ws3.Range("F1").PasteSpecial($xlPasteFormats)

$style1Donor.Copy()
$ws1.Range("M1").PasteSpecial($xlPasteFormats)
$ws2.Range("M1").PasteSpecial($xlPasteFormats)
$ws3.Range("M1").PasteSpecial($xlPasteFormats)
$ws2.Range("C2").PasteSpecial($xlPasteFormats)
$ws2.Range("D2").PasteSpecial($xlPasteFormats)
$ws2.Range("E23").PasteSpecial($xlPasteFormats)
$ws2.Range("E24").PasteSpecial($xlPasteFormats)
$ws2.Range("E25").PasteSpecial($xlPasteFormats)
$ws2.Range("E26").PasteSpecial($xlPasteFormats)
$ws2.Range("E27").PasteSpecial($xlPasteFormats)

# --------------------------------------------------------------
# 3. Add the new "Pt4" sheet after "Pt3".
# --------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "Pt4"

# --------------------------------------------------------------
# 4. Populate Pt4's data.
# --------------------------------------------------------------
$ws4.Range("A1").Value = "Ieeg name"
$ws4.Range("B1").Value = "Current"
$ws4.Range("C1").Value = "Main stim start time"
$ws4.Range("D1").Value = "Main stim end time"
$ws4.Range("E1").Value = "Electrodes"
$ws4.Range("F1").Value = "Clinical Stim seizure elecs"
$ws4.Range("G1").Value = "Suspected SOZ anatomic"
$ws4.Range("H1").Value = "Afterdischarges"
$ws4.Range("I1").Value = "Clinical effects"
$ws4.Range("J1").Value = "Seizures"
$ws4.Range("K1").Value = "Current test electrodes"
$ws4.Range("L1").Value = "Electrode"
$ws4.Range("M1").Value = "Anatomical target"
$ws4.Range("N1").Value = "Other"
$ws4.Range("B2").Value = 3
$ws4.Range("E2").Value = "LM3"
$ws4.Range("F2").Value = "LJ1"
$ws4.Range("H2").Value = "LB1"
$ws4.Range("I2").Value = "LA1 and LA2: weird cephalic feeling"
$ws4.Range("K2").Value = "LM3"
$ws4.Range("L2").Value = "LA"
$ws4.Range("M2").Value = "left amygdala"
$ws4.Range("N2").Value = "Looks like nice cceps LC6->LB9; LA1->LC1;LA2->LL1;LJ3->LH5"
$ws4.Range("E3").Value = "LD1"
$ws4.Range("F3").Value = "LH3"
$ws4.Range("H3").Value = "LB2"
$ws4.Range("K3").Value = "LD1"
$ws4.Range("L3").Value = "LB"
$ws4.Range("M3").Value = "left anterior hippocampus"
$ws4.Range("E4").Value = "LD2"
$ws4.Range("H4").Value = "LA1"
$ws4.Range("L4").Value = "LC"
$ws4.Range("M4").Value = "left posterior hippocampus"
$ws4.Range("E5").Value = "LD3"
$ws4.Range("H5").Value = "LA2"
$ws4.Range("L5").Value = "LD"
$ws4.Range("M5").Value = "left temporal pole"
$ws4.Range("E6").Value = "LD8"
$ws4.Range("H6").Value = "LA3"
$ws4.Range("L6").Value = "LE"
$ws4.Range("M6").Value = "left anterior cingulate"
$ws4.Range("E7").Value = "LD9"
$ws4.Range("H7").Value = "LH5"
$ws4.Range("L7").Value = "LF"
$ws4.Range("M7").Value = "left mid cingulate"
$ws4.Range("E8").Value = "LL1"
$ws4.Range("L8").Value = "LG"
$ws4.Range("M8").Value = "left parietal MEG dipole"
$ws4.Range("E9").Value = "LL2"
$ws4.Range("L9").Value = "LH"
$ws4.Range("M9").Value = "left SMA"
$ws4.Range("E10").Value = "LL3"
$ws4.Range("L10").Value = "LI"
$ws4.Range("M10").Value = "left frontal eye field"
$ws4.Range("E11").Value = "LC5"
$ws4.Range("L11").Value = "LJ"
$ws4.Range("M11").Value = "left superior frontal gyrus"
$ws4.Range("E12").Value = "LC6"
$ws4.Range("L12").Value = "LK"
$ws4.Range("M12").Value = "left frontal pole"
$ws4.Range("E13").Value = "LC7"
$ws4.Range("L13").Value = "LL"
$ws4.Range("M13").Value = "left orbitofrontal gyrus"
$ws4.Range("E14").Value = "LC8"
$ws4.Range("L14").Value = "LM"
$ws4.Range("M14").Value = "left orbitofrontal"
$ws4.Range("E15").Value = "LB1"
$ws4.Range("E16").Value = "LB2"
$ws4.Range("E17").Value = "LB3"
$ws4.Range("E18").Value = "LB9"
$ws4.Range("E19").Value = "LB10"
$ws4.Range("E20").Value = "LA1"
$ws4.Range("E21").Value = "LA2"
$ws4.Range("E22").Value = "LA3"
$ws4.Range("E23").Value = "LA8"
$ws4.Range("E24").Value = "LA9"
$ws4.Range("E25").Value = "LJ4"
$ws4.Range("E26").Value = "LJ3"
$ws4.Range("E27").Value = "LJ2"
$ws4.Range("E28").Value = "LH5"
$ws4.Range("E29").Value = "LH4"
$ws4.Range("E30").Value = "LJ1"
$ws4.Range("E31").Value = "LH3"
$ws4.Range("E32").Value = "LH2"
$ws4.Range("E33").Value = "LH1"

# --------------------------------------------------------------
# 5. Apply the right font style to every Pt4 cell, matching the
#    pattern used on the other three sheets (header row + most
#    data cells use the "Arial/theme-1" style; the long "E" column
#    electrode list, and the F1 header, use the "default" style).
# --------------------------------------------------------------
$style2Cells = @("F1", "E4", "E5", "E6", "E7", "E8", "E9", "E10", "E11", "E12", "E13", "E14", "E15", "E16", "E17", "E18", "E19", "E20", "E21", "E22", "E23", "E24", "E25", "E26", "E27", "E28", "E29", "E30", "E31")
$style1Cells = @("A1", "B1", "C1", "D1", "E1", "G1", "H1", "I1", "J1", "K1", "L1", "M1", "N1", "B2", "E2", "F2", "H2", "I2", "K2", "L2", "M2", "N2", "E3", "F3", "H3", "K3", "L3", "M3", "H4", "L4", "M4", "H5", "L5", "M5", "H6", "L6", "M6", "H7", "L7", "M7", "L8", "M8", "L9", "M9", "L10", "M10", "L11", "M11", "L12", "M12", "L13", "M13", "L14", "M14", "E32", "E33")

$style2Donor.Copy()
foreach ($addr in $style2Cells) {
    $ws4.Range($addr).PasteSpecial($xlPasteFormats)
}

$style1Donor.Copy()
foreach ($addr in $style1Cells) {
    $ws4.Range($addr).PasteSpecial($xlPasteFormats)
}
